# Update stats for 2026-01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (A1:H1) loses its special style (bold font, thin border,
# centered/top alignment) and reverts to the default "Normal" style.
$ws.Range("A1:H1").Style = "Normal"

# Update the last data row (month 2026-01, serial date 46023) with the
# refreshed statistics.
$ws.Range("B26").Value = 6526
$ws.Range("C26").Value = 1015
$ws.Range("D26").Value = 6074807
$ws.Range("E26").Value = 930.8622433343548
$ws.Range("F26").Value = 10.29237789420314
$ws.Range("G26").Value = 7.749469214437377
$ws.Range("H26").Value = 26.50806043341476
